$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 12 (pushes existing row 12 "The air is..." and below down by one)
$ws.Rows.Item(12).Insert()
$ws.Rows.Item(12).RowHeight = 12

# New sentence about sunrise, using a distinct font entry (Arial 10 black) + wrap text
$ws.Range("D12").Value = "The sun rose at `$sunrise_other`$ this morning!"
$ws.Range("D12").Font.ThemeFont = 1
$ws.Range("D12").WrapText = $true

# Narrow column E per the diff
$ws.Columns.Item(5).ColumnWidth = 4.85546875

# Update sheet view: scroll/freeze pane position and active cell selection
$ws.Application.ActiveWindow.ScrollColumn = 3
$sheetView = $ws.Application.ActiveWindow
$ws.Range("D14").Select()
